$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 87: finish filling in the entry that only had Date/Start Time so far.
# Set the Interruption (D) minutes before the Stop Time (C) so the shared
# formula in E87 picks up the correct interruption value when it recalculates.
$ws.Range("D87").Value = 20
$ws.Range("C87").Value = 0.87847222222222221
$ws.Range("F87").Value = "Coding"

# Row 88: new time-log entry.
$ws.Range("A88").Value = 41929
$ws.Range("B88").Value = 0.71250000000000002
$ws.Range("D88").Value = 30
$ws.Range("C88").Value = 0.84583333333333333
$ws.Range("F88").Value = "Coding"

# Update the window scroll position / active selection to match where the
# user ended up after entering the new rows.
$win = $excel.Windows.Item(1)
$win.ScrollRow = 70
$win.ScrollColumn = 1
$ws.Range("A89").Select()
